# Generate Report for Handoff
#
# The localization status report is updated because "b.md" has become
# ready to be handed off again (a newer source revision exists upstream),
# so its status flips from "Handed back: in sync with en-US" to
# "Ready for handoff" on every sheet, its handoff artifacts / datetimes are
# refreshed, the "Content Duplicate" flag drops to False and an Error Detail
# message is recorded explaining that the handback file is stale.

$wb = $excel.ActiveWorkbook

# Helper: assign a literal TEXT value to a cell without letting the host
# auto-coerce look-alike literals (e.g. "False"/"True") into native Boolean
# cells. Routing the text through a temporary formula and then collapsing it
# to a static value via Copy/PasteSpecial(values) keeps the stored cell type
# as a plain (shared) string, matching how this workbook was authored.
function Set-TextValue {
    param($range, [string]$text)

    $escaped = $text.Replace("""", """""")
    $range.Formula = '="' + $escaped + '"'
    $range.Copy() | Out-Null
    $range.PasteSpecial(-4163) | Out-Null
    $excel.CutCopyMode = $false
}

# ---------------------------------------------------------------------
# Sheet "Overview": refresh the summary row for b.md
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
Set-TextValue $overview.Range("E3") "Ready for handoff"
Set-TextValue $overview.Range("F3") "Ready for handoff"
Set-TextValue $overview.Range("G3") "2016-08-26 16:38:27"

# ---------------------------------------------------------------------
# Sheet "zh-cn": refresh the b.md detail row (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
Set-TextValue $zhcn.Range("C3") "Ready for handoff"
Set-TextValue $zhcn.Range("F3") "False"
Set-TextValue $zhcn.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
Set-TextValue $zhcn.Range("H3") "2016-08-26 16:38:23"
Set-TextValue $zhcn.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fbb4a64c4a1651420aef10b05ed20b2256b38f8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de2641f33ddb51b0eeb8daf0955b4cd65c6b016c/e2e/b.md."

# Widen the "Error Detail" column (P / column 16) so the new message fits.
$zhcn.Columns.Item(16).ColumnWidth = 39.1666666666667

# ---------------------------------------------------------------------
# Sheet "de-de": refresh the b.md detail row (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
Set-TextValue $dede.Range("C3") "Ready for handoff"
Set-TextValue $dede.Range("F3") "False"
Set-TextValue $dede.Range("G3") "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
Set-TextValue $dede.Range("H3") "2016-08-26 16:38:27"
Set-TextValue $dede.Range("P3") "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1fbb4a64c4a1651420aef10b05ed20b2256b38f8/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/de2641f33ddb51b0eeb8daf0955b4cd65c6b016c/e2e/b.md."

# Widen the "Error Detail" column (P / column 16) so the new message fits.
$dede.Columns.Item(16).ColumnWidth = 39.1666666666667
